$d = $word.ActiveDocument

# Locate the paragraph ending in "If so, explain why you need to do
# some adjustments." -- the new answer paragraph goes right after it.
$anchor = $d.Content
$found = $anchor.Find.Execute("If so, explain why you need to do some adjustments.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "anchor sentence not found"
}

$anchor.Collapse(0)
$anchorPara = $anchor.Paragraphs(1)

# Create a fresh paragraph right after the anchor paragraph.
$anchor.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newRange = $newPara.Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:iCs/><w:lang w:val="en-CA"/></w:rPr>'

$run1 = '<w:r>' + $rPr + '<w:t>With the new model, I needed to separate my game into three separate classes: the Model, View, and Controller.</w:t></w:r>'
$run2 = '<w:r>' + $rPr + '<w:t xml:space="preserve"> I needed to do this adjustment because it will make the future development a lot easier for me by following the </w:t></w:r>'
$run3 = '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>MVC design pattern. With all three of those components in a single class, as my current development was heading towards, it would be a lot harder to keep my code clean.</w:t></w:r>'
$run4 = '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>'

$pPr = '<w:pPr><w:widowControl/><w:spacing w:before="120" w:after="120"/><w:jc w:val="both"/>' + $rPr + '</w:pPr>'

$xml = '<w:p ' + $ns + '>' + $pPr + $run1 + $run2 + $run3 + $run4 + '</w:p>'

[void]$newRange.InsertXML($xml)
